$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the regression coefficient text values per the diff:
# -2.821*** -> -2.82***  (C Lag, column B)
# -0.012*   -> -0.01*    (A Lag, column C)
# -0.467*** -> -0.47***  (C Lag, column C)
$ws.Range("B3").Value = "-2.82***"
$ws.Range("C2").Value = "-0.01*"
$ws.Range("C3").Value = "-0.47***"
